$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("emissiondata")

# Add new scenario "H2 balanced": a duplicate of the "H2 heavy" row (row 5)
# with the same Year/emission/group/price, but a new Scenario label.
$ws.Range("A6").Value = "H2 balanced"
$ws.Range("B6").Value = 2035
$ws.Range("C6").Value = "CO2"
$ws.Range("D6").Value = "ETS_CO2"
$ws.Range("E6").Value = 105

# Match formatting of the "H2 heavy" scenario label cell (A5): wrap text.
$ws.Range("A6").WrapText = $true

# Move the selection to A7, as left by the author after adding the new row.
$ws.Range("A7").Select() | Out-Null
